$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column map: B=2, C=3, E=5, F=6
# row => (B, C, E, F)  ($null means "leave unchanged / not set")
$data = @{
    3  = @(0.6880700722480534, $null, 0.6397674343333133, $null)
    4  = @(0.9973951713340502, $null, 0.01962905978144387, $null)
    5  = @(68.62777675996794, $null, 1.255803321474957, $null)
    6  = @(0.469963, 0.530037, 0.9981314999999999, 0.0018685)
    7  = @(26.505977, 26.9429704, 32.508251, 16.0651154)
    8  = @(17, 21, 26, 9)
    9  = @(21, 24, 27, 10)
    10 = @(24, 24, 30, 13)
    11 = @(24, 25, 30, 13)
    12 = @(24, 25, 31, 14)
    13 = @(24, 26, 31, 14)
    14 = @(24, 27, 31, 16)
    15 = @(24, 27, 31, 16)
    16 = @(24, 27, 33, 17)
    17 = @(24, 27, 33, 17)
    18 = @(27, 27, 34, 17)
    19 = @(28, 28, 34, 17)
    20 = @(31, 28, 34, 17)
    21 = @(31, 28, 34, 18)
    22 = @(31, 28, 34, 19)
    23 = @(31, 29, 34, 20)
    24 = @(31, 30, 34, 20)
    25 = @(31, 30, 38, 20)
    26 = @(34, 31, 41, 20)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $b = $vals[0]
    $c = $vals[1]
    $e = $vals[2]
    $f = $vals[3]

    if ($null -ne $b) { $ws.Cells.Item($row, 2).Value = $b }
    if ($null -ne $c) { $ws.Cells.Item($row, 3).Value = $c }
    if ($null -ne $e) { $ws.Cells.Item($row, 5).Value = $e }
    if ($null -ne $f) { $ws.Cells.Item($row, 6).Value = $f }
}
